# RPA datasets push 2024-06-20
# Applies the IPO tracker update to the "02_38커뮤니케이션(최근일자기준)" sheet:
#  - Removes the "이노그리드" row (old row 12), shifting subsequent rows up.
#  - Finalizes 이노스페이스's confirmed offer price (확정공모가) from "-" to 43300.
#  - Appends a brand-new IPO entry (미래에셋비전스팩5호) as the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# 1) Drop the 이노그리드 row entirely; every row below shifts up by one.
$ws.Rows("12:12").Delete()

# 2) 이노스페이스 (now row 14) got its confirmed offering price finalized.
#    Column D in this table stores values as text even when numeric-looking
#    (e.g. "18000", "7000"), so force a text number format before writing the
#    value, then clear the format again so no stray style sticks to the cell.
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "43300"
$ws.Range("D14").ClearFormats()

# 3) Append the new IPO entry as row 21.
$ws.Range("A21").Value = "미래에셋비전스팩5호"
$ws.Range("B21").Value = "2024.06.03~06.04"
$ws.Range("C21").Value = "2,000~2,000"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2000"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = 9500
$ws.Range("F21").Value = "미래에셋증권"
